# ReporteAsistencia.xlsx — trim the report down:
#   - retitle the header banner (B2)
#   - drop the "Asistido" column (E) from the table header
#   - drop the sample data row entirely
#   - pull the header row up to close the resulting gap

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the banner title text in the merged title cell (B2:E4)
$ws.Range("B2").Value = 'Reporte de Asistentes tema: ""'

# Bring the header row's formatting up from row 6 to row 5 (columns A-D only;
# column E / "Asistido" is being dropped)
$ws.Range("A6:D6").Copy()
$ws.Range("A5:D5").PasteSpecial(-4122)

# Move the header labels themselves up a row, skipping column E ("Asistido")
$ws.Range("A5").Value = $ws.Range("A6").Value()
$ws.Range("B5").Value = $ws.Range("B6").Value()
$ws.Range("C5").Value = $ws.Range("C6").Value()
$ws.Range("D5").Value = $ws.Range("D6").Value()

# Remove the now-duplicate header row (old row 6) and the sample data row
# beneath it (old row 7); each Delete shifts the rows below up by one.
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
